$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A (shifts existing columns B:F right by one).
$ws.Columns("A").Insert()

# New header + matching column width for the inserted "Employee_ID" column.
$ws.Range("A1").Value = "Employee_ID"
$ws.Range("A1").Font.Bold = $true
$ws.Columns("A").ColumnWidth = 12

# New data value. Format as text first so the leading zeros in "003" are
# preserved instead of Excel coercing it to the number 3.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "003"
